$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial of 45181 for every data row (2-300).
# Update it to 45182 for each row.
for ($r = 2; $r -le 300; $r++) {
    $ws.Cells.Item($r, 3).Value = 45182
}
